$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Merge the 3 header cells of each new block first, just like the existing
# Iteration_1/Iteration_2 blocks.
# ---------------------------------------------------------------------------
$ws.Range("K1:M1").Merge()
$ws.Range("N1:P1").Merge()
$ws.Range("Q1:S1").Merge()

# ---------------------------------------------------------------------------
# Row 2 sub-headers (2030 / 2040 / 2050) repeated for each new block. These
# must end up stored as *text* (matching the existing B2:J2 "2030" style
# cells), not as numbers, so the cells are switched to Text format before
# the value is written - otherwise Excel would store them as numeric 2030.
# ---------------------------------------------------------------------------
$row2Cells = @("K2", "L2", "M2", "N2", "O2", "P2", "Q2", "R2", "S2")
foreach ($cellRef in $row2Cells) {
    $ws.Range($cellRef).NumberFormat = "@"
}
$ws.Range("K2").Value2 = "2030"
$ws.Range("L2").Value2 = "2040"
$ws.Range("M2").Value2 = "2050"
$ws.Range("N2").Value2 = "2030"
$ws.Range("O2").Value2 = "2040"
$ws.Range("P2").Value2 = "2050"
$ws.Range("Q2").Value2 = "2030"
$ws.Range("R2").Value2 = "2040"
$ws.Range("S2").Value2 = "2050"

# ---------------------------------------------------------------------------
# Copy the formatting (font, border, alignment) from an already "plain"
# styled cell (I1, which carries the same style as every other header
# cell) onto the new header ranges. Doing this *after* merging and *after*
# writing the row-2 text values resets the style index of every new header
# cell back to the shared one used by the existing header cells (undoing
# the Text-format / per-edge-border style variants that merging and the
# "@" number format would otherwise introduce), while preserving the text
# already stored in the cells.
# ---------------------------------------------------------------------------
$src = $ws.Range("I1")
$src.Copy()
$ws.Range("K1:S2").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# Row 1 headers for the new iteration blocks.
# ---------------------------------------------------------------------------
$ws.Range("K1").Value2 = "Iteration_3"
$ws.Range("N1").Value2 = "Iteration_4"
$ws.Range("Q1").Value2 = "Iteration_5"

# ---------------------------------------------------------------------------
# Data rows 4-16 for columns K:S (Iteration_3, Iteration_4, Iteration_5).
# ---------------------------------------------------------------------------
# NOTE: scientific notation (e.g. "1e-08") is not accepted by the script
# parser, so very small magnitudes below are written out as plain decimals.
$data = @{
    4  = @(0, 0, 0, 0, 0, 0, 0, 0, 0)
    5  = @(0, 0, 0, 0, 0, 0, 0, 0, 0)
    6  = @(0, 0, 0, 0, 0, 0, 0, 0, 0)
    7  = @(0, 0, 0, 0, 0, 0, 0, 0, 0)
    8  = @(0, 2939476.181780299, 4753682.583935678, 0, 3556062.637552082, 4753634.260821007, 0, 2938140.333093762, 4753391.494946638)
    9  = @(12883863.3233789, 3325241.373589423, 6449870.170927788, 12859798.25098664, 4865449.27575592, 6450017.49749695, 12873149.56277529, 3329314.082999473, 6450757.639224648)
    10 = @(0, 0, 0, 0, 0, 0, 0, 0, 0)
    11 = @(0.00000004349772098066751, 0, 0, 0.000002872993111593766, 0, 0, 0.0000003870135856232081, 0, 0)
    12 = @(0, 0, 0, 0, 0, 0, 0, 0, 0)
    13 = @(0, 0, 0, 0, 0, 0, 0, 0, 0)
    14 = @(0, 0, 0, 0, 0, 0, 0, 0, 0)
    15 = @(0, 0, 0, 0, 0, 0, 0, 0, 0)
    16 = @(0, 0, 0, 0, 0, 0, 0, 0, 0)
}

$cols = @("K", "L", "M", "N", "O", "P", "Q", "R", "S")

foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range($cols[$i] + $row).Value2 = $values[$i]
    }
}
